$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 134, shifting existing rows 134-136 down to 135-137
$ws.Rows.Item(134).Insert()

# Copy formatting/style of column D from the row above (row 133) into new row's D134 (date style)
$ws.Range("D133").Copy()
$ws.Range("D134").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new row 134 data
$ws.Range("A134").Value = 3
$ws.Range("B134").Value = "Femacal de La Calera"
$ws.Range("C134").Value = "Coquimbo"
$ws.Range("D134").Value = 44448
$ws.Range("E134").Value = 5
$ws.Range("F134").Value = 100112001
$ws.Range("G134").Value = "Berenjena"
$ws.Range("H134").Value = "Sin especificar"
$ws.Range("I134").Value = "Primera"
$ws.Range("J134").Value = 130
$ws.Range("K134").Value = 9500
$ws.Range("L134").Value = 10000
$ws.Range("M134").Value = 9692
$ws.Range("N134").Value = "$/caja 60 unidades"
$ws.Range("O134").Value = "Región de Arica y Parinacota"
$ws.Range("P134").Value = 162
$ws.Range("Q134").Value = 60
$ws.Range("R134").Value = "Hortaliza"
